# Hortaliza, Vega Modelo de Temuco - Coliflor
# Commit: "Fruta / hortaliza, semanal"
#
# Insert two new weekly observation rows above the current row 592,
# pushing the existing data (old rows 592:666) down to 594:668.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 592 - everything from old row 592
# onward shifts down by two rows (old 592 -> new 594, ..., old 666 -> new 668).
$ws.Range("A592:A593").EntireRow.Insert()

# New row 592: Región Metropolitana, 2023-07-17
$ws.Range("A592").Value = 10
$ws.Range("B592").Value = "Vega Modelo de Temuco"
$ws.Range("C592").Value = "La Araucanía"
$ws.Range("D592").Value = 45124
$ws.Range("E592").Value = 9
$ws.Range("F592").Value = 100112008
$ws.Range("G592").Value = "Coliflor"
$ws.Range("H592").Value = "Sin especificar"
$ws.Range("I592").Value = "Primera"
$ws.Range("J592").Value = 2800
$ws.Range("K592").Value = 1000
$ws.Range("L592").Value = 1000
$ws.Range("M592").Value = 1000
$ws.Range("N592").Value = "$/unidad"
$ws.Range("O592").Value = "Región Metropolitana"
$ws.Range("P592").Value = 1000
$ws.Range("Q592").Value = 1
$ws.Range("R592").Value = "Hortaliza"

# New row 593: Región del Maule, 2023-07-17
$ws.Range("A593").Value = 10
$ws.Range("B593").Value = "Vega Modelo de Temuco"
$ws.Range("C593").Value = "La Araucanía"
$ws.Range("D593").Value = 45124
$ws.Range("E593").Value = 9
$ws.Range("F593").Value = 100112008
$ws.Range("G593").Value = "Coliflor"
$ws.Range("H593").Value = "Sin especificar"
$ws.Range("I593").Value = "Primera"
$ws.Range("J593").Value = 3000
$ws.Range("K593").Value = 1000
$ws.Range("L593").Value = 1000
$ws.Range("M593").Value = 1000
$ws.Range("N593").Value = "$/unidad"
$ws.Range("O593").Value = "Región del Maule"
$ws.Range("P593").Value = 1000
$ws.Range("Q593").Value = 1
$ws.Range("R593").Value = "Hortaliza"
